$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ B = -152.69155300104572;  C = 0.08247937813401086;  D = 423.854304641 }
    3  = @{ B = -149.80697366092807;  C = 0.05280934461389653;  D = 254.751864734 }
    4  = @{ B = -151.1044281555749;   C = 0.0996243877259803;   D = 362.976633382 }
    5  = @{ B = -150.7026638682956;   C = 0.08533525115220324;  D = 1003.487114047 }
    6  = @{ B = -150.52128312101613;  C = 0.09977938483127843;  D = 263.92695321 }
    7  = @{ B = -149.30562672348597;  C = 0.08055856980811969;  D = 493.548096787 }
    8  = @{ B = -148.661520255107;    C = 0.09939557762595787;  D = 308.777549555 }
    9  = @{ B = -150.18260019047344;  C = 0.09070489220644856;  D = 620.859802629 }
    10 = @{ B = -151.08461268611245;  C = 0.030069965494367006; D = 647.406392567 }
    11 = @{ B = -146.9222127693745;   C = 0.09681741853795965;  D = 722.842084454 }
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row].B
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("D$row").Value = $values[$row].D
}
